# daily auto push: 2025-09-30 07:28 UTC
# Append the new daily data row (row 40) to Sheet1, extending the table
# that currently ends at row 39 (A1:D39 -> A1:D40).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds a date-like label (e.g. "2025/09/30") that must be stored
# as literal text, not auto-converted to a date serial number. Assigning it
# with a leading apostrophe forces Excel to keep it as text; resetting the
# cell style back to "Normal" afterwards drops the quote-prefix formatting
# flag so the cell ends up with no special style, matching the rest of the
# column.
$ws.Range("A40").Value = "'2025/09/30"
$ws.Range("A40").Style = "Normal"

$ws.Range("B40").Value = "火"
$ws.Range("C40").Value = 16
$ws.Range("D40").Value = 152
